$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
  @("ECs", "Fgf1", "Cspg4", "ECs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 1.577167666666667, 4.731503, 0.05487405574265763, 0.05487405574265762, 1.831655761302111, 16.484901851719, 0.006973671370507471, 0.00697367137050747),
  @("ECs", "Fgf1", "Cspg4", "FAPs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 12.088216, 36.264648, 0.4205827019109694, 0.4205827019109694, 14.03874232792267, 126.348680951304, 0.05344976797417882, 0.05344976797417882),
  @("ECs", "Fgf1", "Cspg4", "M2", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 2, 0.6666666666666666, 0.06147399999999999, 0.184422, 0.002138851673172859, 0.002138851673172859, 0.07139330120066667, 0.6425397108059999, 0.0002718160427018071, 0.0002718160427018071),
  @("ECs", "Fgf1", "Cspg4", "sCs", 2, 0.6666666666666666, 1.161357666666667, 3.484073, 0.1270850363824361, 0.1270850363824361, 3, 1, 15.01473333333333, 45.0442, 0.5224043906732001, 0.5224043906732001, 17.43747566962222, 156.9372810266, 0.066389780995048, 0.066389780995048),
  @("FAPs", "Fgf1", "Cspg4", "ECs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 1.577167666666667, 4.731503, 0.05487405574265763, 0.05487405574265762, 8.15107640278478, 73.359687625063, 0.03103363052701016, 0.03103363052701015),
  @("FAPs", "Fgf1", "Cspg4", "FAPs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 12.088216, 36.264648, 0.4205827019109694, 0.4205827019109694, 62.47399960817868, 562.2659964736081, 0.2378575448909317, 0.2378575448909316),
  @("FAPs", "Fgf1", "Cspg4", "M2", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 2, 0.6666666666666666, 0.06147399999999999, 0.184422, 0.002138851673172859, 0.002138851673172859, 0.3177083079846666, 2.859374771862, 0.001209612296357472, 0.001209612296357472),
  @("FAPs", "Fgf1", "Cspg4", "sCs", 3, 1, 5.168173666666667, 15.504521, 0.565542861868062, 0.565542861868062, 3, 1, 15.01473333333333, 45.0442, 0.5224043906732001, 0.5224043906732001, 77.59874942535555, 698.3887448282, 0.2954420741537627, 0.2954420741537627),
  @("sCs", "Fgf1", "Cspg4", "ECs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 1.577167666666667, 4.731503, 0.05487405574265763, 0.05487405574265762, 4.430103630287222, 39.870932672585, 0.01686675384514, 0.01686675384514),
  @("sCs", "Fgf1", "Cspg4", "FAPs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 12.088216, 36.264648, 0.4205827019109694, 0.4205827019109694, 33.95456977537334, 305.59112797836, 0.1292753890458589, 0.1292753890458589),
  @("sCs", "Fgf1", "Cspg4", "M2", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 2, 0.6666666666666666, 0.06147399999999999, 0.184422, 0.002138851673172859, 0.002138851673172859, 0.1726742161433333, 1.55406794529, 0.0006574233341135805, 0.0006574233341135805),
  @("sCs", "Fgf1", "Cspg4", "sCs", 3, 1, 2.808898333333333, 8.426695, 0.3073721017495019, 0.3073721017495019, 3, 1, 15.01473333333333, 45.0442, 0.5224043906732001, 0.5224043906732001, 42.17485943544445, 379.573734919, 0.1605725355243894, 0.1605725355243894),
)

for ($i = 0; $i -lt $data.Count; $i++) {
  $r = $i + 2
  $row = $data[$i]
  for ($j = 0; $j -lt $row.Count; $j++) {
    $c = $j + 1
    $ws.Cells.Item($r, $c).Value = $row[$j]
  }
}

Write-Host "done"